$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '62.037.54'
Set-TextValue $ws.Range("E2") '  -1.93%  '
Set-TextValue $ws.Range("D3") '2.591.42'
Set-TextValue $ws.Range("E3") '  -4.39%  '
Set-TextValue $ws.Range("E4") '  -0.13%  '
Set-TextValue $ws.Range("D5") '553.84'
Set-TextValue $ws.Range("E5") '  -1.14%  '
Set-TextValue $ws.Range("D6") '155.40'
Set-TextValue $ws.Range("E6") '  -0.49%  '
Set-TextValue $ws.Range("E7") '  -0.19%  '
Set-TextValue $ws.Range("E8") '  +0.73%  '
Set-TextValue $ws.Range("E9") '  -2.16%  '
Set-TextValue $ws.Range("E10") '  -3.01%  '
Set-TextValue $ws.Range("D11") '5.49'
Set-TextValue $ws.Range("E11") '  -1.46%  '
Set-TextValue $ws.Range("D12") '0.366'
Set-TextValue $ws.Range("E12") '  -1.31%  '
Set-TextValue $ws.Range("D13") '3.048.01'
Set-TextValue $ws.Range("E13") '  -4.57%  '
Set-TextValue $ws.Range("D14") '25.63'
Set-TextValue $ws.Range("E14") '  -2.64%  '
Set-TextValue $ws.Range("D15") '61.898.21'
Set-TextValue $ws.Range("E15") '  -1.98%  '
Set-TextValue $ws.Range("E16") '  -1.99%  '
Set-TextValue $ws.Range("D17") '2.592.62'
Set-TextValue $ws.Range("E17") '  -4.48%  '
Set-TextValue $ws.Range("D18") '11.67'
Set-TextValue $ws.Range("E18") '  -3.78%  '
Set-TextValue $ws.Range("D19") '4.56'
Set-TextValue $ws.Range("E19") '  -2.17%  '
Set-TextValue $ws.Range("D20") '339.67'
Set-TextValue $ws.Range("E20") '  -3.08%  '
Set-TextValue $ws.Range("D21") '6.06'
Set-TextValue $ws.Range("E21") '  -5.87%  '
Set-TextValue $ws.Range("D22") '0.998'
Set-TextValue $ws.Range("E22") '  -0.04%  '
Set-TextValue $ws.Range("D23") '0.500'
Set-TextValue $ws.Range("E23") '  -2.00%  '
Set-TextValue $ws.Range("D24") '62.57'
Set-TextValue $ws.Range("E24") '  -2.47%  '
Set-TextValue $ws.Range("E25") '  -0.33%  '
Set-TextValue $ws.Range("D26") '1.00'
Set-TextValue $ws.Range("E26") '  -0.05%  '
Set-TextValue $ws.Range("D27") '8.09'
Set-TextValue $ws.Range("E27") '  -1.04%  '
Set-TextValue $ws.Range("D28") '0.0₃0840'
Set-TextValue $ws.Range("E28") '  -5.32%  '
Set-TextValue $ws.Range("D29") '1.93'
Set-TextValue $ws.Range("E29") '  -1.06%  '
Set-TextValue $ws.Range("D30") '7.11'
Set-TextValue $ws.Range("E30") '  -0.49%  '
Set-TextValue $ws.Range("E31") '  -3.22%  '
Set-TextValue $ws.Range("D32") '160.03'
Set-TextValue $ws.Range("E32") '  -3.48%  '
Set-TextValue $ws.Range("E33") '  -0.04%  '
Set-TextValue $ws.Range("B34") 'NEARProtocol'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D34") '4.73'
Set-TextValue $ws.Range("E34") '  -1.54%  '
Set-TextValue $ws.Range("B35") 'EthereumClassic'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D35") '19.27'
Set-TextValue $ws.Range("E35") '  -2.59%  '
Set-TextValue $ws.Range("D36") '1.43'
Set-TextValue $ws.Range("E36") '  -3.15%  '
Set-TextValue $ws.Range("D37") '1.78'
Set-TextValue $ws.Range("E37") '  +0.82%  '
Set-TextValue $ws.Range("D38") '341.40'
Set-TextValue $ws.Range("E38") '  -0.83%  '
Set-TextValue $ws.Range("D39") '6.05'
Set-TextValue $ws.Range("E39") '  -0.91%  '
Set-TextValue $ws.Range("D40") '0.899'
Set-TextValue $ws.Range("E40") '  -6.30%  '
Set-TextValue $ws.Range("D41") '3.93'
Set-TextValue $ws.Range("E41") '  -2.51%  '
Set-TextValue $ws.Range("D42") '37.56'
Set-TextValue $ws.Range("E42") '  -2.26%  '
Set-TextValue $ws.Range("D43") '20.61'
Set-TextValue $ws.Range("E43") '  -3.40%  '
Set-TextValue $ws.Range("E44") '  -0.11%  '
Set-TextValue $ws.Range("B45") 'Mantle'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D45") '0.609'
Set-TextValue $ws.Range("E45") '  -2.11%  '
Set-TextValue $ws.Range("B46") 'Maker'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D46") '2.134.92'
Set-TextValue $ws.Range("E46") '  +1.66%  '
Set-TextValue $ws.Range("E47") '  -4.58%  '
Set-TextValue $ws.Range("B48") 'Hedera'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D48") '0.0550'
Set-TextValue $ws.Range("E48") '  -3.81%  '
Set-TextValue $ws.Range("B49") 'WhiteBITCoin'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D49") '10.94'
Set-TextValue $ws.Range("E49") '  -1.07%  '
Set-TextValue $ws.Range("E50") '  -1.70%  '
Set-TextValue $ws.Range("E51") '  -1.64%  '
